# Updated symbol list on Tue Dec 13 05:39:18 UTC 2022 with GitHub Actions
#
# This script reapplies the latest crypto price/volume refresh to Sheet1.
# Columns: A=Rank, B=Coin, C=Link, D=Price, E=Volume(1h), F=Data, G=Hora
#
# Numeric-looking "Price" values are stored as text in this workbook (the
# column uses inline strings, not numbers), so we force each cell's number
# format to Text ("@") before writing the value - this prevents Excel from
# re-interpreting strings like "267.72" or "0.001636" as floating point
# numbers (which would silently drop meaningful trailing/leading zeros,
# e.g. "0.1120" -> "0.112"). We then restore the cell style to "Normal" so
# we don't leave a stray custom number format behind on cells that didn't
# have one originally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $value) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Price (column D) refreshes ---
Set-TextValue $ws "D2"  "267.72"
Set-TextValue $ws "D3"  "21.55"
Set-TextValue $ws "D4"  "6.256"
Set-TextValue $ws "D5"  "0.06206"
Set-TextValue $ws "D6"  "3.573"
Set-TextValue $ws "D7"  "6.533"
Set-TextValue $ws "D8"  "1.393"
Set-TextValue $ws "D9"  "0.8244"
Set-TextValue $ws "D10" "0.1638"
Set-TextValue $ws "D11" "0.08205"
Set-TextValue $ws "D12" "0.03558"
Set-TextValue $ws "D13" "0.03189"
Set-TextValue $ws "D14" "0.09201"
Set-TextValue $ws "D16" "0.001636"
Set-TextValue $ws "D17" "0.04691"
Set-TextValue $ws "D18" "0.006462"
Set-TextValue $ws "D19" "0.006185"
Set-TextValue $ws "D22" "3.722"
Set-TextValue $ws "D23" "2.247"
Set-TextValue $ws "D24" "0.01361"
Set-TextValue $ws "D25" "0.3318"
Set-TextValue $ws "D28" "0.0002715"
Set-TextValue $ws "D40" "0.04712"
Set-TextValue $ws "D41" "0.006971"

# --- Rows 42/43 swapped places (BKEXToken <-> CEJI) with new data ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.004003"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1120"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Remaining price refreshes ---
Set-TextValue $ws "D44" "0.01180"
Set-TextValue $ws "D45" "0.00006336"
Set-TextValue $ws "D46" "0.0009906"
Set-TextValue $ws "D47" "0.00000000751"

Set-TextValue $ws "D48" "0.9807"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"

Set-TextValue $ws "D49" "0.002305"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"

Set-TextValue $ws "D50" "0.00001901"
Set-TextValue $ws "D51" "0.01241"
